# Add a new worksheet "phoneme pronunciation" with example phonetic
# pronunciation data, matching the CSV import behaviour described in the
# commit message ("Add Phoneme_Pronunciation.csv").

$wb = $excel.ActiveWorkbook

# Add the new sheet, then move it to the end of the workbook (after the
# current last sheet, "manual translations").
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "phoneme pronunciation"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-acquire a live reference to the sheet after the move.
$newSheet = $wb.Worksheets.Item("phoneme pronunciation")

# Header row
$newSheet.Range("A1").Value = "Text"
$newSheet.Range("B1").Value = "Phonetic Pronunciation"
$newSheet.Range("C1").Value = "Case Sensitive (True/False)"
$newSheet.Range("D1").Value = "Phonetic Alphabet"

# Example data row
$newSheet.Range("A2").Value = "ThioJoe"
$newSheet.Range("B2").Value = "ˈθioʊd͡ʒoʊ"
$newSheet.Range("C2").Value = $false
$newSheet.Range("D2").Value = "ipa"

# Make this new sheet the active / selected tab, and select cell D6 as the
# worksheet's saved selection (matches the committed workbook's state).
$newSheet.Select()
$newSheet.Range("D6").Select()
